# Update "want to go" counts (column F) and "lowest price" (column G)
# across the four worksheets, per the regenerated site data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 368
$ws.Range("F4").Value = 417
$ws.Range("F5").Value = 1142
$ws.Range("F8").Value = 971
$ws.Range("F9").Value = 1626
$ws.Range("F10").Value = 6107
$ws.Range("F12").Value = 1763
$ws.Range("F13").Value = 451
$ws.Range("F14").Value = 6008
$ws.Range("F18").Value = 97
$ws.Range("F19").Value = 1661
$ws.Range("F22").Value = 148
$ws.Range("F23").Value = 1425
$ws.Range("F25").Value = 260
$ws.Range("F28").Value = 32

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 312
$ws.Range("G4").Value = 168
$ws.Range("F8").Value = 384

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9505
$ws.Range("F3").Value = 2236
$ws.Range("F4").Value = 632
$ws.Range("F5").Value = 199

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9505
$ws.Range("F3").Value = 2236
$ws.Range("F4").Value = 632
$ws.Range("F5").Value = 368
$ws.Range("F6").Value = 417
$ws.Range("F7").Value = 1142
$ws.Range("F11").Value = 312
$ws.Range("G11").Value = 168
$ws.Range("F12").Value = 971
$ws.Range("F13").Value = 199
$ws.Range("F14").Value = 1626
$ws.Range("F15").Value = 6107
$ws.Range("F17").Value = 1763
$ws.Range("F20").Value = 451
$ws.Range("F23").Value = 6008
$ws.Range("F27").Value = 97
$ws.Range("F28").Value = 1661
$ws.Range("F31").Value = 148
$ws.Range("F32").Value = 1425
$ws.Range("F35").Value = 260
